$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new LP numbers / statuses
$ws.Range("A2").Value = "LP-047185"
$ws.Range("B2").Value = "Possui linhas de compra e apontamento!"

$ws.Range("A3").Value = "LP-048980"
$ws.Range("B3").Value = "Compromisso pendente!"

$ws.Range("A4").Value = "LP-049043"
$ws.Range("B4").Value = "Compromisso pendente!"

# Remove the now-unused last row (row 5)
$ws.Rows.Item(5).Delete()

# Refresh the active selection to match the new data extent
[void]$ws.Range("A2:B4").Select()
